$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new D (price) value; $null means D unchanged
$priceUpdates = @{
    2  = "36.888.22"
    3  = "1.977.60"
    5  = "245.09"
    7  = "60.76"
    10 = "0.0800"
    12 = "15.00"
    13 = "0.844"
    14 = "22.06"
    15 = "2.270.93"
    16 = "5.47"
    17 = "1.977.51"
    18 = "36.803.01"
    19 = "70.28"
    20 = "0.0`u{2083}0861"
    22 = "229.64"
    26 = "0.147"
    27 = "9.29"
    28 = "163.01"
    29 = "19.52"
    30 = "1.36"
    31 = "0.121"
    32 = "4.87"
    34 = "4.52"
    39 = "5.55"
    40 = "0.0998"
    44 = "16.44"
    45 = "1.369.93"
    46 = "90.04"
    48 = "7.27"
}

# Map of row -> new E (volume) value
$volumeUpdates = @{
    2  = "  +0.74%  "
    3  = "  +0.70%  "
    4  = "  +0.05%  "
    5  = "  +0.18%  "
    6  = "  +1.38%  "
    7  = "  +2.37%  "
    8  = "  +0.02%  "
    9  = "  +1.83%  "
    10 = "  -1.78%  "
    11 = "  +0.67%  "
    12 = "  +9.32%  "
    13 = "  +1.66%  "
    14 = "  -1.09%  "
    15 = "  +0.81%  "
    16 = "  +3.75%  "
    17 = "  +0.56%  "
    18 = "  +0.66%  "
    19 = "  +0.31%  "
    20 = "  +0.27%  "
    21 = "  +1.95%  "
    22 = "  +0.16%  "
    23 = "  +0.01%  "
    24 = "  +2.10%  "
    25 = "  +0.41%  "
    26 = "  +4.19%  "
    27 = "  +0.38%  "
    28 = "  +1.85%  "
    29 = "  +0.59%  "
    30 = "  +17.64%  "
    31 = "  +1.46%  "
    32 = "  +3.30%  "
    34 = "  +5.67%  "
    35 = "  -0.03%  "
    36 = "  -0.49%  "
    37 = "  -1.11%  "
    38 = "  +0.01%  "
    39 = "  -7.71%  "
    40 = "  +1.03%  "
    41 = "  +0.75%  "
    42 = "  +0.51%  "
    43 = "  +0.57%  "
    44 = "  +1.74%  "
    45 = "  +0.78%  "
    46 = "  +2.44%  "
    47 = "  -0.27%  "
    48 = "  +1.71%  "
    49 = "  -0.74%  "
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Value = $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}

# Rows 50 and 51 swap coin identity (B, C) and get new D/E values
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "46.25"
$ws.Range("E50").Value = "  +5.52%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.99"
$ws.Range("E51").Value = "  +11.76%  "
